# Add an "Electrode Locations" column (C) derived from the file name in
# column A, and re-sort the data rows (2..last) by that electrode location
# (letter prefix, then numeric suffix) from A1 up to O15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row based on column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
$firstDataRow = 2

# Read existing data rows (filename + value) into an array of objects,
# extracting the electrode location (e.g. "A11") from the filename.
$records = @()
for ($r = $firstDataRow; $r -le $lastRow; $r++) {
    $fname = $ws.Cells.Item($r, 1).Value2
    $value = $ws.Cells.Item($r, 2).Value2

    if ($fname -match '^([A-Z]+)(\d+)_') {
        $letter = $matches[1]
        $num = [int]$matches[2]
        $loc = "$letter$num"
    }
    else {
        $letter = ""
        $num = 0
        $loc = ""
    }

    $records += [PSCustomObject]@{
        FileName = $fname
        Value    = $value
        Letter   = $letter
        Num      = $num
        Loc      = $loc
    }
}

# Sort by electrode location: letter prefix alphabetically, then numeric
# suffix ascending (so A2, A3, A8 ... A14, B6, B13 ... matches "A1-O15" order).
$sorted = $records | Sort-Object Letter, Num

# Write the header for the new column, matching the style of the existing
# header cells (A1 / B1).
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 3).Value = "Electrode Locations"

# Write the sorted rows back out, filling in column C with the electrode
# location extracted from the file name.
for ($i = 0; $i -lt $sorted.Count; $i++) {
    $r = $firstDataRow + $i
    $rec = $sorted[$i]
    $ws.Cells.Item($r, 1).Value = $rec.FileName
    $ws.Cells.Item($r, 2).Value = $rec.Value
    $ws.Cells.Item($r, 3).Value = $rec.Loc
}
